$d = $word.ActiveDocument

# The document's headers contain the "BTec_Logo-Orange" picture (currently
# stored/exported as image2.jpg) which must be renamed to image1.jpg, and the
# footers contain the Pearson logo picture (currently image1.png) which must
# be renamed to image2.png. InlineShape has no settable Name in the Word
# object model, so each picture is temporarily converted to a floating Shape
# (which does expose a settable Name), renamed, and converted back to an
# InlineShape so the drawing remains <wp:inline> exactly as before.

foreach ($sec in $d.Sections) {

    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($ishp in $hdr.Range.InlineShapes) {
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp = $ishp.ConvertToShape()
                    $shp.Name = "image1.jpg"
                    $null = $shp.ConvertToInlineShape()
                }
            }
        }
    }

    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($ishp in $ftr.Range.InlineShapes) {
                if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp = $ishp.ConvertToShape()
                    $shp.Name = "image2.png"
                    $null = $shp.ConvertToInlineShape()
                }
            }
        }
    }
}
